$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing schedule times back by 3 hours (20:00->17:00, 21:00->18:00, 22:00->19:00)
$ws.Range("A2").Value = 0.70833333333333337
$ws.Range("A3").Value = 0.75
$ws.Range("A4").Value = 0.79166666666666663

# Add two new rows to the schedule, reusing the time slots freed up above
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 0.83333333333333337

$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 0.875

$ws.Range("B5").Value = "Higurashi : When They Cry"
$ws.Range("C5").Value = "HR_Bumper.mp4"
$ws.Range("D5").Value = "HR.txt"
$ws.Range("E5").Value = "HR_p.txt"

$ws.Range("B6").Value = "Dragon Ball Super"
$ws.Range("C6").Value = "DBS_Bumper.mp4"
$ws.Range("D6").Value = "DBS.txt"
$ws.Range("E6").Value = "DBS_p.txt"

$ws.Range("B6").Select()
